$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing empty placeholder rows at the bottom of the sheet
$ws.Rows.Item(1048575).Delete()
$ws.Rows.Item(1048575).Delete()

# Update SeedMass units from "g" to "mg" (row 120, column E)
$ws.Cells.Item(120, 5).Value = "mg"

# Insert two new rows for dispersal kernel parameters after SeedLongevity (row 121)
$ws.Rows.Item(122).Insert()
$ws.Rows.Item(123).Insert()

# Row 122: DispersalDistance
$ws.Cells.Item(122, 1).Value = "DispersalDistance"
$ws.Cells.Item(122, 2).Value = "Regeneration"
$ws.Cells.Item(122, 3).Value = "Distance parameter for dispersal kernel"
$ws.Cells.Item(122, 4).Value = "Numeric"
$ws.Cells.Item(122, 5).Value = "m"

# Row 123: DispersalShape
$ws.Cells.Item(123, 1).Value = "DispersalShape"
$ws.Cells.Item(123, 2).Value = "Regeneration"
$ws.Cells.Item(123, 3).Value = "Shape parameter for dispersal kernel"
$ws.Cells.Item(123, 4).Value = "Numeric"

# Restore view state: gridlines visible, selection/scroll near the edited area
$ws.Range("E120").Select()
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 1
